{"js": "// 1. Update the date heading paragraph (first paragraph of the body).\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[0].insertText(\"2023-11-22 Wednesday\", \"Replace\");\nawait context.sync();\n\n// 2. Update each answer cell in the table by (row, column) position,\n//    which avoids any ambiguity from values that repeat elsewhere.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\ntable.getCell(0, 0).value = \"34\u00d772=2448\"; // was \"30\u00d737=1110\"\ntable.getCell(0, 1).value = \"31\u00d783=2573\"; // was \"13\u00d783=1079\"\ntable.getCell(0, 2).value = \"74\u00d744=3256\"; // was \"96\u00d723=2208\"\ntable.getCell(0, 3).value = \"44\u00d769=3036\"; // was \"29\u00d768=1972\"\ntable.getCell(0, 4).value = \"92\u00d731=2852\"; // was \"61\u00d753=3233\"\ntable.getCell(4, 0).value = \"69\u00d739=2691\"; // was \"90\u00d767=6030\"\ntable.getCell(4, 1).value = \"68\u00d742=2856\"; // was \"45\u00d735=1575\"\ntable.getCell(4, 2).value = \"34\u00d774=2516\"; // was \"85\u00d718=1530\"\ntable.getCell(4, 3).value = \"48\u00d720=960\"; // was \"48\u00d732=1536\"\ntable.getCell(4, 4).value = \"30\u00d749=1470\"; // was \"86\u00d771=6106\"\ntable.getCell(9, 0).value = \"66\u00d736=2376\"; // was \"86\u00d741=3526\"\ntable.getCell(9, 1).value = \"16\u00d762=992\"; // was \"59\u00d785=5015\"\ntable.getCell(9, 2).value = \"16\u00d757=912\"; // was \"54\u00d776=4104\"\ntable.getCell(9, 3).value = \"41\u00d767=2747\"; // was \"50\u00d787=4350\"\ntable.getCell(9, 4).value = \"25\u00d797=2425\"; // was \"17\u00d795=1615\"\ntable.getCell(14, 0).value = \"20\u00d737=740\"; // was \"34\u00d772=2448\"\ntable.getCell(14, 1).value = \"26\u00d746=1196\"; // was \"45\u00d713=585\"\ntable.getCell(14, 2).value = \"62\u00d778=4836\"; // was \"39\u00d742=1638\"\ntable.getCell(14, 3).value = \"39\u00d713=507\"; // was \"78\u00d745=3510\"\ntable.getCell(14, 4).value = \"38\u00d757=2166\"; // was \"33\u00d798=3234\"\ntable.getCell(19, 0).value = \"84\u00d768=5712\"; // was \"26\u00d756=1456\"\ntable.getCell(19, 1).value = \"20\u00d792=1840\"; // was \"14\u00d715=210\"\ntable.getCell(19, 2).value = \"53\u00d781=4293\"; // was \"71\u00d725=1775\"\ntable.getCell(19, 3).value = \"85\u00d737=3145\"; // was \"69\u00d721=1449\"\ntable.getCell(19, 4).value = \"19\u00d723=437\"; // was \"16\u00d784=1344\"\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the date heading (first paragraph of the document).\n$d.Paragraphs.Item(1).Range.Text = \"2023-11-22 Wednesday\"\n\n# 2. Update each answer cell in the table by (row, column) position,\n#    which avoids any ambiguity from values that repeat elsewhere.\n$tbl = $d.Tables.Item(1)\n\n$tbl.Cell(1,1).Range.Text = \"34\u00d772=2448\" # was \"30\u00d737=1110\"\n$tbl.Cell(1,2).Range.Text = \"31\u00d783=2573\" # was \"13\u00d783=1079\"\n$tbl.Cell(1,3).Range.Text = \"74\u00d744=3256\" # was \"96\u00d723=2208\"\n$tbl.Cell(1,4).Range.Text = \"44\u00d769=3036\" # was \"29\u00d768=1972\"\n$tbl.Cell(1,5).Range.Text = \"92\u00d731=2852\" # was \"61\u00d753=3233\"\n$tbl.Cell(5,1).Range.Text = \"69\u00d739=2691\" # was \"90\u00d767=6030\"\n$tbl.Cell(5,2).Range.Text = \"68\u00d742=2856\" # was \"45\u00d735=1575\"\n$tbl.Cell(5,3).Range.Text = \"34\u00d774=2516\" # was \"85\u00d718=1530\"\n$tbl.Cell(5,4).Range.Text = \"48\u00d720=960\" # was \"48\u00d732=1536\"\n$tbl.Cell(5,5).Range.Text = \"30\u00d749=1470\" # was \"86\u00d771=6106\"\n$tbl.Cell(10,1).Range.Text = \"66\u00d736=2376\" # was \"86\u00d741=3526\"\n$tbl.Cell(10,2).Range.Text = \"16\u00d762=992\" # was \"59\u00d785=5015\"\n$tbl.Cell(10,3).Range.Text = \"16\u00d757=912\" # was \"54\u00d776=4104\"\n$tbl.Cell(10,4).Range.Text = \"41\u00d767=2747\" # was \"50\u00d787=4350\"\n$tbl.Cell(10,5).Range.Text = \"25\u00d797=2425\" # was \"17\u00d795=1615\"\n$tbl.Cell(15,1).Range.Text = \"20\u00d737=740\" # was \"34\u00d772=2448\"\n$tbl.Cell(15,2).Range.Text = \"26\u00d746=1196\" # was \"45\u00d713=585\"\n$tbl.Cell(15,3).Range.Text = \"62\u00d778=4836\" # was \"39\u00d742=1638\"\n$tbl.Cell(15,4).Range.Text = \"39\u00d713=507\" # was \"78\u00d745=3510\"\n$tbl.Cell(15,5).Range.Text = \"38\u00d757=2166\" # was \"33\u00d798=3234\"\n$tbl.Cell(20,1).Range.Text = \"84\u00d768=5712\" # was \"26\u00d756=1456\"\n$tbl.Cell(20,2).Range.Text = \"20\u00d792=1840\" # was \"14\u00d715=210\"\n$tbl.Cell(20,3).Range.Text = \"53\u00d781=4293\" # was \"71\u00d725=1775\"\n$tbl.Cell(20,4).Range.Text = \"85\u00d737=3145\" # was \"69\u00d721=1449\"\n$tbl.Cell(20,5).Range.Text = \"19\u00d723=437\" # was \"16\u00d784=1344\"\n"}
